# Update the Overview sheet's yearly income-statement table:
#  - drop the oldest year (1396/12) and its "published" date,
#  - shift the remaining four years one column to the left (D<-E, E<-F, F<-G, G<-H),
#  - add the newly reported year (1401/12) in column H with fresh figures,
#  - the previously blank/"-" cell D15 becomes a real 0.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$cols = @("D", "E", "F", "G", "H")

function Set-RowValues($row, $values) {
    for ($i = 0; $i -lt 5; $i++) {
        $addr = $cols[$i] + $row
        $ws.Range($addr).Value = $values[$i]
    }
}

# Row 8: period headers (financial period labels)
Set-RowValues 8 @("12 ماهه منتهی به 1397/12", "12 ماهه منتهی به 1398/12", "12 ماهه منتهی به 1399/12", "12 ماهه منتهی به 1400/12", "12 ماهه منتهی به 1401/12")

# Row 9: publish dates
Set-RowValues 9 @("1399-02-01 (9)", "1400-02-01 (8)", "1401-02-07 (9)", "1402-02-06 (9)", "1402-02-06 (2)")

# Row 11: فروش (Sales)
Set-RowValues 11 @(651742, 907909, 1839050, 3456701, 5516295)

# Row 12: بهای تمام شده کالای فروش رفته (COGS)
Set-RowValues 12 @(-361369, -515874, -773456, -1620409, -2070101)

# Row 13: سود (زیان) ناخالص (Gross profit)
Set-RowValues 13 @(290373, 392035, 1065594, 1836292, 3446194)

# Row 14: هزینه های عمومی, اداری و تشکیلاتی
Set-RowValues 14 @(-52767, -69154, -96522, -149888, -231795)

# Row 15: هزینه کاهش ارزش دریافتنی‌ها (D was the "-" placeholder, now numeric 0)
Set-RowValues 15 @(0, 0, 0, 0, 0)

# Row 16: خالص سایر درامدها (هزینه ها) ی عملیاتی
Set-RowValues 16 @(-7590, -1343, 629, 2196, -15654)

# Row 17: سود (زیان) عملیاتی
Set-RowValues 17 @(230016, 321538, 969701, 1688600, 3198745)

# Row 18: هزینه های مالی
Set-RowValues 18 @(-55164, -12740, -35660, -50660, -62657)

# Row 19: خالص سایر درامدها و هزینه های غیرعملیاتی
Set-RowValues 19 @(40050, 13778, 41558, 83365, 84016)

# Row 20: سود (زیان) خالص عملیات در حال تداوم قبل از مالیات
Set-RowValues 20 @(214902, 322576, 975599, 1721305, 3220104)

# Row 21: مالیات
Set-RowValues 21 @(-63878, -76832, -216918, -313135, -458085)

# Row 22: سود (زیان) خالص عملیات در حال تداوم
Set-RowValues 22 @(151024, 245744, 758681, 1408170, 2762019)

# Row 23: سود (زیان) عملیات متوقف شده پس از اثر مالیاتی
Set-RowValues 23 @(0, 0, 0, 0, 0)

# Row 24: سود (زیان) خالص
Set-RowValues 24 @(151024, 245744, 758681, 1408170, 2762019)

# Row 25: سود هر سهم پس از کسر مالیات
Set-RowValues 25 @(755, 1229, 3793, 4694, 9207)

# Row 26: سرمایه
Set-RowValues 26 @(200000, 200000, 200000, 300000, 300000)

# Row 27: سود هر سهم بر اساس آخرین سرمایه
Set-RowValues 27 @(503, 819, 2529, 4694, 9207)
